$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update StartingInventories (C) and SetupCosts (E) ---
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 0.01701
$ws.Range("E3").Value = 0.012564
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0.012627
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 0.0124875
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0.0139995
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 0.0144
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0.00144
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 0.002232
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 0.00216
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 0.00153
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 0.001494
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 0.002142
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 0.0004635
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 0.012609
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 0.0004815
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 0.001494
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 0.0005715
$ws.Range("C19").Value = 0
$ws.Range("E19").Value = 0.016173
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 0.000396
$ws.Range("E21").Value = 0.006581249999999999
$ws.Range("E22").Value = 0.00663975
$ws.Range("E23").Value = 0.006599249999999999
$ws.Range("E24").Value = 0.00678375
$ws.Range("E25").Value = 0.00732375
$ws.Range("E26").Value = 0.007281000000000001
$ws.Range("E27").Value = 0.00755775
$ws.Range("E28").Value = 0.0076455

# --- Capacity sheet: update capacity values (B) ---
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 80
$ws.Range("B4").Value = 30
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 30
$ws.Range("B7").Value = 20
$ws.Range("B8").Value = 80
$ws.Range("B9").Value = 160
$ws.Range("B10").Value = 40
$ws.Range("B11").Value = 120
$ws.Range("B12").Value = 80
$ws.Range("B13").Value = 80
$ws.Range("B14").Value = 50
$ws.Range("B15").Value = 40
$ws.Range("B16").Value = 30
$ws.Range("B17").Value = 200
$ws.Range("B18").Value = 30
$ws.Range("B19").Value = 100
$ws.Range("B20").Value = 10
$ws.Range("B21").Value = 10
$ws.Range("B22").Value = 5
$ws.Range("B23").Value = 15
$ws.Range("B24").Value = 20
$ws.Range("B25").Value = 10
$ws.Range("B26").Value = 10
$ws.Range("B27").Value = 20
$ws.Range("B28").Value = 25

# --- ProcessingTime sheet: update diagonal processing times ---
$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("E5").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("G7").Value = 2
$ws.Range("I9").Value = 4
$ws.Range("J10").Value = 1
$ws.Range("K11").Value = 3
$ws.Range("L12").Value = 2
$ws.Range("N14").Value = 5
$ws.Range("O15").Value = 2
$ws.Range("P16").Value = 3
$ws.Range("Q17").Value = 5
$ws.Range("R18").Value = 3
$ws.Range("S19").Value = 5
$ws.Range("T20").Value = 1
$ws.Range("U21").Value = 2
$ws.Range("V22").Value = 1
$ws.Range("W23").Value = 3
$ws.Range("X24").Value = 4
$ws.Range("Y25").Value = 2
$ws.Range("Z26").Value = 2
$ws.Range("AA27").Value = 4
